$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fuels_by_sectors")

# Give all remaining "16_others" subfuel rows (16_02 .. 16_x_ammonia, and 16_x_efuel)
# to the Buildings sector by marking column D ("Buildings") with "a",
# matching the already-marked 16_01_biogas (row 54) and 16_x_hydrogen (row 64) rows.
$rows = 55,56,57,58,59,60,61,62,63,65
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "a"
}

# Reflect the resulting active selection on the sheet (last edited cell D65)
$ws.Range("D65").Select()
